$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2 through 89 is updated from
# serial date 45177 (2023-09-08) to serial date 45178 (2023-09-09).
for ($r = 2; $r -le 89; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
